$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.272.12"
$ws.Range("E2").Value = "  +3.18%  "

# Row 3
$ws.Range("D3").Value = "2.268.91"
$ws.Range("E3").Value = "  +2.55%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.84"
$ws.Range("E5").Value = "  +3.79%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.01"
$ws.Range("E6").Value = "  +3.63%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.589"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  +0.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.73"
$ws.Range("E10").Value = "  +4.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  +2.01%  "

# Row 12
$ws.Range("E12").Value = "  +2.94%  "

# Row 13
$ws.Range("E13").Value = "  -0.38%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.614.35"
$ws.Range("E14").Value = "  +2.57%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.880"
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.55"
$ws.Range("E16").Value = "  +2.01%  "

# Row 17
$ws.Range("D17").Value = "2.268.50"
$ws.Range("E17").Value = "  +2.54%  "

# Row 18
$ws.Range("D18").Value = "44.176.30"
$ws.Range("E18").Value = "  +3.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.28"
$ws.Range("E19").Value = "  -4.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +4.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("E21").Value = "  +1.09%  "

# Row 22
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.26"
$ws.Range("E22").Value = "  +4.37%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.48"
$ws.Range("E23").Value = "  +1.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.72"
$ws.Range("E24").Value = "  +0.25%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").Value = "  +1.94%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.30"
$ws.Range("E27").Value = "  +2.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.58"
$ws.Range("E28").Value = "  +11.98%  "

# Row 29
$ws.Range("E29").Value = "  -1.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.48"
$ws.Range("E30").Value = "  +2.23%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "162.72"
$ws.Range("E31").Value = "  +4.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.58"
$ws.Range("E32").Value = "  -0.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0881"
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.73"
$ws.Range("E34").Value = "  -2.45%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.118"
$ws.Range("E35").Value = "  +12.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.99"
$ws.Range("E36").Value = "  +1.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.12"
$ws.Range("E37").Value = "  +2.44%  "

# Row 38
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.92"
$ws.Range("E39").Value = "  +4.11%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.38"
$ws.Range("E40").Value = "  -0.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.82"
$ws.Range("E41").Value = "  +30.30%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0326"
$ws.Range("E42").Value = "  -0.12%  "

# Row 43
$ws.Range("E43").Value = "  +0.14%  "

# Row 44
$ws.Range("D44").Value = "1.789.24"
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.206"
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "84.79"
$ws.Range("E46").Value = "  -2.96%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.40"
$ws.Range("E47").Value = "  +1.36%  "

# Row 48
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "77.07"
$ws.Range("E48").Value = "  +0.71%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.75"
$ws.Range("E49").Value = "  +6.76%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.39"
$ws.Range("E50").Value = "  -0.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "105.96"
$ws.Range("E51").Value = "  +2.66%  "
